# Atualização de bases das ligas, do dia: 06-04-2024 às 15:39
#
# This script applies the league-odds update to the "Bolivia Primera
# División" workbook's single worksheet:
#   1. Rows 148-150 (match ids 146-148) get their betting-odds columns
#      rotated/updated to the newly scraped values.
#   2. Rows 210-212 (match ids 208-210) get updated closing-odds values.
#   3. A brand-new row 213 (match id 211, fixture 8039392) is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $Values) {
    foreach ($col in $Values.Keys) {
        $ws.Cells.Item($Row, $col).Value = $Values[$col]
    }
}

# Column index reference:
# A=1 B=2 C=3 D=4 E=5 F=6 G=7 H=8 I=9 J=10 K=11 L=12 M=13 N=14 O=15 P=16
# Q=17 R=18 S=19 T=20 U=21 V=22 W=23 X=24 Y=25 Z=26 AA=27 AB=28 AC=29

# ---------------------------------------------------------------------
# Row 148 (id 146)
# ---------------------------------------------------------------------
Set-Row 148 @{
    2  = 7532420
    6  = "Club Aurora"
    7  = "Vaca Diez"
    11 = 1.333
    12 = 5
    13 = 8
    14 = 1.3
    15 = 6.5
    16 = 7
    17 = -1.5
    18 = 1.8
    19 = 2
    20 = 3.25
    23 = 0.3
    26 = 0.8
    28 = -0.5
    29 = 0.425
}

# ---------------------------------------------------------------------
# Row 149 (id 147)
# ---------------------------------------------------------------------
Set-Row 149 @{
    2  = 7532419
    6  = "Oriente Petrolero"
    7  = "Jorge Wilstermann"
    8  = 3
    11 = 2.2
    12 = 2.5
    13 = 4.5
    14 = 2.375
    15 = 2.45
    16 = 4.5
    17 = -0.25
    18 = 1.9
    19 = 1.9
    20 = 2
    23 = 1.375
    26 = 0.8999999999999999
    28 = 0.95
    29 = -1
}

# ---------------------------------------------------------------------
# Row 150 (id 148)
# ---------------------------------------------------------------------
Set-Row 150 @{
    2  = 7532421
    6  = "Guabira"
    7  = "Independiente Petrolero"
    8  = 2
    11 = 1.4
    12 = 4.5
    13 = 7.5
    14 = 1.333
    15 = 5.5
    16 = 9.5
    17 = -1.5
    18 = 1.85
    19 = 1.95
    20 = 3
    23 = 0.333
    26 = 0.8500000000000001
    28 = -1
    29 = 0.9750000000000001
}

# ---------------------------------------------------------------------
# Row 210 (id 208) - odds refresh only
# ---------------------------------------------------------------------
Set-Row 210 @{
    14 = 1.666
    15 = 3.8
    16 = 4.75
    18 = 1.8
    19 = 2
    20 = 3
    21 = 1.9
    22 = 1.9
}

# ---------------------------------------------------------------------
# Row 211 (id 209) - odds refresh only
# ---------------------------------------------------------------------
Set-Row 211 @{
    14 = 3
    17 = 0
    18 = 2.05
    19 = 1.75
    21 = 1.9
    22 = 1.9
}

# ---------------------------------------------------------------------
# Row 212 (id 210) - odds refresh only
# ---------------------------------------------------------------------
Set-Row 212 @{
    14 = 1.95
    15 = 3.5
    16 = 3.8
    17 = -0.5
    18 = 1.975
    19 = 1.825
    20 = 2.25
    21 = 1.825
    22 = 1.975
}

# ---------------------------------------------------------------------
# Row 213 (id 211) - brand-new fixture row
# ---------------------------------------------------------------------
# Copy formatting from the row above (A: bold/bordered id style, E: date
# number format) so the new row matches the existing table styling.
$ws.Cells.Item(212, 1).Copy($ws.Cells.Item(213, 1))
$ws.Cells.Item(212, 5).Copy($ws.Cells.Item(213, 5))

Set-Row 213 @{
    1  = 211
    2  = 8039392
    3  = "Bolivia Primera División"
    4  = "Bolivia Apertura"
    5  = 45389.70833333334
    6  = "Oriente Petrolero"
    7  = "Jorge Wilstermann"
    11 = 2
    12 = 3.25
    13 = 3.4
    14 = 1.75
    15 = 3.4
    16 = 4
    17 = -0.5
    18 = 1.8
    19 = 2
    20 = 2.25
    21 = 1.775
    22 = 2.025
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 0
}
